# edit.ps1 - applies the "fixing technical variations section" commit.
#
# Summary of changes (see diff):
#  1. Heading "4. Replicating ..." -> "4. Results of Replicating ..." (+ bookmark rename)
#  2. Remove " [1]" citation marker in the "sample NMR spectral data" sentence.
#  3. Replace trailing "[1]" citation in "target image" sentence with a proper
#     "(Xia J 2009)" citation, split across 4 runs like elsewhere in the doc.
#  4. Heading "5. Change" -> "5. Technical Variations" (+ bookmark rename)
#  5. Heading "5.1 Change in imputation method" -> "5.1 Imputation method" (+ bookmark rename)
#  6. Reword "The first of the changes ..." sentence.
#  7. Reword "... detection limit." sentence (merge two sentences with ", and").
#  8. Reword final sentence of 5.1 paragraph, removing a sentence and changing tense.
#  9. Insert new paragraph after 5.1: "As we expected, ... did not change."
# 10. Heading "5.2 Change in filtering method" -> "5.2 Filtering method" (+ bookmark rename)
# 11. Insert new paragraph after 5.2: "The output remained constant ..."
# 12. Heading "5.3 Change in Component for Variable Importance" -> "5.3 Component for Variable Importance" (+ bookmark rename)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: rename a bookmark while keeping it anchored to the same range.
# ---------------------------------------------------------------------------
function Rename-Bookmark($oldName, $newName) {
    $bm = $d.Bookmarks($oldName)
    $bmRange = $d.Range($bm.Start, $bm.End)
    $d.Bookmarks.Add($newName, $bmRange)
    $d.Bookmarks($oldName).Delete()
}

# ---------------------------------------------------------------------------
# Helper: plain whole-text replace (single resulting run).
# ---------------------------------------------------------------------------
function Replace-Text($oldText, $newText) {
    $ok = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) { throw ("Replace-Text: not found: " + $oldText) }
}

# ---------------------------------------------------------------------------
# Helper: insert a zero-length bookmark at an absolute document position and
# immediately delete it again. Because of how the runtime serializes runs,
# doing this forces a run split at that exact position (the bookmark start/
# end tags act as a seam) without altering any visible text.
# ---------------------------------------------------------------------------
function Add-Seam($pos) {
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add("__seam__", $r)
    $d.Bookmarks("__seam__").Delete()
}

# ---------------------------------------------------------------------------
# Helper: replace a whole paragraph/run's text with a new string that should
# be split into several runs. $pieces is an array of substrings of $newFull
# (in order, concatenation == $newFull) each of which will become its own
# <w:r> in the saved document.
# ---------------------------------------------------------------------------
function Replace-TextMultiRun($oldText, $pieces) {
    $newFull = [string]::Join("", $pieces)
    $rng = $d.Content
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw ("Replace-TextMultiRun: not found: " + $oldText) }
    $paraStart = $rng.Start
    $rng.Text = $newFull
    $paraEnd = $rng.Start + $newFull.Length

    # Compute seam boundary offsets (document positions) between pieces.
    $pos = $paraStart
    $boundaries = @()
    for ($i = 0; $i -lt ($pieces.Length - 1); $i++) {
        $pos = $pos + $pieces[$i].Length
        $boundaries += $pos
    }

    foreach ($b in $boundaries) {
        Add-Seam $b
    }
}

# ===========================================================================
# 1. Heading 4
# ===========================================================================
Replace-Text "4. Replicating PLSDA Analysis of NMR Spectral Bin Data" "4. Results of Replicating PLSDA Analysis of NMR Spectral Bin Data"
Rename-Bookmark "X7134d7f5f0a977c366712219bb72618c5ca115e" "X913da3038ee0c9a266e5bc5e2f95fef9b8b766d"

# ===========================================================================
# 2. Remove "[1]" citation marker
# ===========================================================================
Replace-Text "We use the sample NMR spectral data that comes with the tutorial [1] in the online tool for this analysis." "We use the sample NMR spectral data that comes with the tutorial in the online tool for this analysis."

# ===========================================================================
# 3. Replace "[1]" citation with "(Xia J 2009)" citation, split into 4 runs
# ===========================================================================
Replace-TextMultiRun "This is the target image that we are trying to replicate. We get this image from [1]." @(
    "This is the target image that we are trying to replicate. We get this image from",
    " ",
    "(Xia J 2009)",
    "."
)

# ===========================================================================
# 4. Heading 5
# ===========================================================================
Replace-Text "5. Change" "5. Technical Variations"
Rename-Bookmark "change" "technical-variations"

# ===========================================================================
# 5. Heading 5.1
# ===========================================================================
Replace-Text "5.1 Change in imputation method" "5.1 Imputation method"
Rename-Bookmark "change-in-imputation-method" "imputation-method"

# ===========================================================================
# 6. Reword first sentence of 5.1 paragraph
# ===========================================================================
Replace-Text "The first of the changes that we are going to implement is to change the imputation method for columns with less than 50% of missing values, which is set by default to" "The first of the changes that we implemented was the imputation method for columns with less than 50% of missing values, which is set by default to"

# ===========================================================================
# 7. Merge two sentences in 5.1 paragraph
# ===========================================================================
Replace-Text ". This method assigns half of the column’s minimum positive value to the missing observations. The assumption of this approach is that most missing values occurred because the levels of abundance metabolites are below the detection limit. This assumption is generalized for all the possible datasets used in the MetaboAnalyst website. We decided to change the imputation method to" ". This method assigns half of the column’s minimum positive value to the missing observations. The assumption of this approach is that most missing values occurred because the levels of abundance metabolites are below the detection limit, and this assumption is generalized for all the possible datasets used in the MetaboAnalyst website. We decided to change the imputation method to"

# ===========================================================================
# 8. Reword final sentence of 5.1 paragraph (drop one sentence, change tense)
# ===========================================================================
Replace-Text ", a more widely used method that imputes the mean value of the column to the missing values. This change shouldn’t impact the output in a significant manner. Due to the nature of the dataset used, data provided for the testing of the tool, we believe that it should not make a difference." ", a more widely used method that imputes the mean value of the column to the missing values. Due to the nature of the dataset used, data provided for the testing of the tool, we believed that it should not make a difference."

# ===========================================================================
# 9. Insert new "As we expected..." paragraph after the 5.1 paragraph
# ===========================================================================
$rng = $d.Content
$found = $rng.Find.Execute("believed that it should not make a difference.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "could not find end of 5.1 paragraph" }
$rng.InsertParagraphAfter()
$newPara = $rng.Paragraphs(1).Next()
$newRange = $newPara.Range
$newRange.MoveEnd(1, -1)
$newRange.InsertAfter("As we expected, because of a lack of missing values in the data, the result did not change.")
$newPara.Style = "Body Text"

# ===========================================================================
# 10. Heading 5.2
# ===========================================================================
Replace-Text "5.2 Change in filtering method" "5.2 Filtering method"
Rename-Bookmark "change-in-filtering-method" "filtering-method"

# ===========================================================================
# 11. Insert new "The output remained constant..." paragraph after 5.2
# ===========================================================================
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("mantain near-constant values throughout the entire experiment.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "could not find end of 5.2 paragraph" }
$rng2.InsertParagraphAfter()
$newPara2 = $rng2.Paragraphs(1).Next()
$newRange2 = $newPara2.Range
$newRange2.MoveEnd(1, -1)
$newRange2.InsertAfter("The output remained constant after changing the filtering method.")
$newPara2.Style = "Body Text"

# ===========================================================================
# 12. Heading 5.3
# ===========================================================================
Replace-Text "5.3 Change in Component for Variable Importance" "5.3 Component for Variable Importance"
Rename-Bookmark "Xb4cc110fdc7b35fe584b201b60e723d67add69d" "component-for-variable-importance"

Write-Output "All edits applied."
